$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells: force text format to avoid numeric auto-conversion
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D10", "D12", "D14", "D17", "D18", "D19", "D21", "D23", "D24", "D26", "D27", "D30", "D31", "D32", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D46", "D47", "D48", "D50", "D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Set new values
$ws.Range("D2").Value = '43.434.74'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '2.271.24'
$ws.Range("E3").Value = '  +1.09%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '266.02'
$ws.Range("E5").Value = '  -2.17%  '
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '93.57'
$ws.Range("E6").Value = '  -4.85%  '
$ws.Range("D7").Value = '0.619'
$ws.Range("E7").Value = '  -1.15%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.603'
$ws.Range("E9").Value = '  -3.94%  '
$ws.Range("D10").Value = '44.11'
$ws.Range("E10").Value = '  -8.46%  '
$ws.Range("E11").Value = '  -1.08%  '
$ws.Range("D12").Value = '7.67'
$ws.Range("E12").Value = '  -6.88%  '
$ws.Range("E13").Value = '  +0.49%  '
$ws.Range("D14").Value = '2.612.14'
$ws.Range("E14").Value = '  +1.16%  '
$ws.Range("E15").Value = '  -1.90%  '
$ws.Range("E16").Value = '  +0.89%  '
$ws.Range("D17").Value = '2.279.89'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").Value = '43.435.72'
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("D19").Value = '0.0000106'
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("E20").Value = '  -1.29%  '
$ws.Range("D21").Value = '71.61'
$ws.Range("E21").Value = '  +0.96%  '
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("D23").Value = '232.95'
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("D24").Value = '8.88'
$ws.Range("E24").Value = '  -8.30%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = '2.48'
$ws.Range("E26").Value = '  -0.83%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '11.20'
$ws.Range("E27").Value = '  -1.79%  '
$ws.Range("E28").Value = '  -2.23%  '
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").Value = '38.67'
$ws.Range("E30").Value = '  -2.99%  '
$ws.Range("D31").Value = '175.37'
$ws.Range("E31").Value = '  +1.21%  '
$ws.Range("D32").Value = '21.71'
$ws.Range("E32").Value = '  +2.92%  '
$ws.Range("E33").Value = '  -3.74%  '
$ws.Range("D34").Value = '5.29'
$ws.Range("E34").Value = '  -6.16%  '
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.0354'
$ws.Range("E36").Value = '  +0.86%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '0.106'
$ws.Range("E37").Value = '  -5.48%  '
$ws.Range("D38").Value = '4.36'
$ws.Range("E38").Value = '  -0.72%  '
$ws.Range("D39").Value = '3.27'
$ws.Range("E39").Value = '  -8.17%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.236'
$ws.Range("E40").Value = '  -6.90%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").Value = '2.34'
$ws.Range("E41").Value = '  +6.70%  '
$ws.Range("E42").Value = '  +14.93%  '
$ws.Range("D43").Value = '11.75'
$ws.Range("E43").Value = '  -6.42%  '
$ws.Range("D44").Value = '61.81'
$ws.Range("E44").Value = '  -0.51%  '
$ws.Range("E45").Value = '  +3.22%  '
$ws.Range("D46").Value = '5.17'
$ws.Range("E46").Value = '  -5.03%  '
$ws.Range("D47").Value = '0.101'
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("D48").Value = '97.01'
$ws.Range("E48").Value = '  -3.51%  '
$ws.Range("E49").Value = '  -0.94%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.491.53'
$ws.Range("E50").Value = '  +1.13%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = '1.47'
$ws.Range("E51").Value = '  +3.63%  '

# Restore default style on price cells (remove the temporary text-format style)
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
